$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core Dataset refactor: rename columns, drop "_duplicates" column ---

# Delete the trailing "_duplicates" column (M). All of its cells (header
# + boolean flag values) disappear and every column to its right would
# shift left - but since it is the last column, nothing shifts; the
# sheet's dimension simply shrinks from A1:M13 to A1:L13.
$ws.Columns("M").Delete()

# Rename header cells to reflect the new Dataset-based column naming
# scheme (PIDN_link/DCDate_link -> PIDN_x/DCDate_x, and the internal
# merge/diff columns get an "_mp_" prefix).
$ws.Range("A1").Value = "PIDN_x"
$ws.Range("B1").Value = "DCDate_x"
$ws.Range("J1").Value = "_mp_merge"
$ws.Range("K1").Value = "_mp_diff_days"
$ws.Range("L1").Value = "_mp_abs_diff_days"

# Adjust column widths to match the new best-fit values for the renamed
# headers (the engine quantizes ColumnWidth to the nearest 1/6 of a
# character, so these are the closest achievable values to the target
# best-fit widths).
$ws.Columns("A").ColumnWidth = 5.833333333333333
$ws.Columns("J").ColumnWidth = 9.5
$ws.Columns("K").ColumnWidth = 11.166666666666666
$ws.Columns("L").ColumnWidth = 14.5
